$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 126, shifting rows 126:130 down to 127:131
$ws.Rows.Item(126).Insert()

# Populate the new row 126 with the new data entry
$ws.Cells.Item(126, 1).Value = 10
$ws.Cells.Item(126, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(126, 3).Value = "La Araucanía"
$ws.Cells.Item(126, 4).Value = 44931
$ws.Cells.Item(126, 5).Value = 9
$ws.Cells.Item(126, 6).Value = "Fruta"
$ws.Cells.Item(126, 7).Value = 100101
$ws.Cells.Item(126, 8).Value = "Berries"
$ws.Cells.Item(126, 9).Value = 100101001
$ws.Cells.Item(126, 10).Value = "Arándano (blue)"
$ws.Cells.Item(126, 11).Value = "Sin especificar"
$ws.Cells.Item(126, 12).Value = "Primera"
$ws.Cells.Item(126, 13).Value = 250
$ws.Cells.Item(126, 14).Value = 1800
$ws.Cells.Item(126, 15).Value = 1800
$ws.Cells.Item(126, 16).Value = 1800
$ws.Cells.Item(126, 17).Value = "$/kilo"
$ws.Cells.Item(126, 18).Value = "Región del Maule"
$ws.Cells.Item(126, 19).Value = 1800
$ws.Cells.Item(126, 20).Value = 1
